$d = $word.ActiveDocument

$d.Content.Find.Execute("{percent}%{/retained}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "{percent}{/retained}", 2)
